# Update "丽水-漫展信息.xlsx" - 展览 (sheet 1) and 全部类型 (sheet 4) sheets:
# Row 2 is overwritten with the data that used to live in row 3, and the
# (now duplicate) row 3 is removed.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("B2").NumberFormat = "@"
    $ws.Range("B2").Value = "2024-06-01"
    $ws.Range("C2").Value = "丽水·动漫游戏展"
    $ws.Range("D2").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Range("E2").Value = "2024.06.01 10:00-06.01 17:00"
    $ws.Range("F2").Value = 388
    $ws.Range("G2").Value = 45
    $ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=84450"
    $ws.Range("I2").Value = "//i2.hdslb.com/bfs/openplatform/202404/tdhb9QSW1713333412467.jpeg"

    # Remove the now-duplicated third row entirely (shifts rows up, nothing below to shift).
    $ws.Rows.Item(3).Delete()
}
